# edit.ps1 - Applies the ORDENANZA N 618 formatting update:
#  - Adds keepNext + spacing (before/after) to the first five paragraphs
#  - Bolds the title lines ("ORDENANZA N 618" and "EL CONCEJO DELIBERANTE...")
#  - Centers/indents the "EL CONCEJO..." paragraph
#  - Underlines the "ARTICULO PRIMERO:" / "ARTICULO SEGUNDO:" labels while
#    keeping the following space in its own (non-underlined) run
#  - Adds a document footer (empty, "Piedepgina" style) -> footer1.xml +
#    footerReference
#  - Triggers creation of footnotes.xml / endnotes.xml (separator parts)
#  - Restarts page numbering at this section starting at 613

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Paragraph 1: "Yerba Buena, 10 de Marzo de 1994"
# ---------------------------------------------------------------------
$p1 = $d.Paragraphs.Item(1)
$p1.Format.KeepWithNext = $true
$p1.Format.SpaceAfter = 12

# ---------------------------------------------------------------------
# Paragraph 2: "ORDENANZA N 618"
# ---------------------------------------------------------------------
$p2 = $d.Paragraphs.Item(2)
$p2.Format.KeepWithNext = $true
$p2.Format.SpaceBefore = 12
$p2.Format.SpaceAfter = 18
$p2.Range.Font.Bold = $true

# ---------------------------------------------------------------------
# Paragraph 3: "EL CONCEJO DELIBERANTE SANCIONA CON FUERZA DE ORDENANZA"
# ---------------------------------------------------------------------
$p3 = $d.Paragraphs.Item(3)
$p3.Format.KeepWithNext = $true
$p3.Format.SpaceBefore = 18
$p3.Format.SpaceAfter = 18
$p3.Format.LeftIndent = 99.2
$p3.Format.RightIndent = 99.2
$p3.Range.Font.Bold = $true

# ---------------------------------------------------------------------
# Paragraph 4: "ARTICULO PRIMERO: ..."
# ---------------------------------------------------------------------
$p4 = $d.Paragraphs.Item(4)
$p4.Format.KeepWithNext = $true
$p4.Format.SpaceAfter = 6

$rng4 = $d.Content
$rng4.Find.ClearFormatting()
$null = $rng4.Find.Execute("ARTICULO PRIMERO: ", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$label4 = $rng4.Duplicate
$label4.End = $label4.End - 1
$label4.Font.Underline = 1

# ---------------------------------------------------------------------
# Paragraph 5: "ARTICULO SEGUNDO: ..."
# ---------------------------------------------------------------------
$p5 = $d.Paragraphs.Item(5)
$p5.Format.KeepWithNext = $true
$p5.Format.SpaceAfter = 6

$rng5 = $d.Content
$rng5.Find.ClearFormatting()
$null = $rng5.Find.Execute("ARTICULO SEGUNDO: ", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$label5 = $rng5.Duplicate
$label5.End = $label5.End - 1
$label5.Font.Underline = 1

# ---------------------------------------------------------------------
# Section / footer / page numbering
# ---------------------------------------------------------------------
$section = $d.Sections.Item(1)

$footer = $section.Footers.Item(1)
$footerPara = $footer.Range.Paragraphs.Item(1)
$footerPara.Style = "Piedepgina"
$footerPara.Range.Font.Name = "Book Antiqua"
$footerPara.Range.Font.Size = 10
$footerPara.Range.Font.Color = 8421504

$footerStyle = $d.Styles.Item("Piedepgina")
$footerStyle.NameLocal = "footer"
$footerStyle.UnhideWhenUsed = $true

# Create footnotes.xml / endnotes.xml (separator parts) without leaving any
# visible footnote reference behind.
$fnRange = $d.Range(0, 0)
$tempNote = $d.Footnotes.Add($fnRange, "", "x")
$tempNote.Delete()

$pageNumbers = $footer.PageNumbers
$pageNumbers.RestartNumberingAtSection = $true
$pageNumbers.StartingNumber = 613
